$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 800
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 1500
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 1500
$ws.Range("M41").Value = 340
$ws.Range("N41").Value = -2380
$ws.Range("H42").Value = 1481.9
$ws.Range("I42").Value = 1260.7142
$ws.Range("J42").Value = 1998
$ws.Range("K42").Value = 3782.1426
$ws.Range("L42").Value = 5994
$ws.Range("M42").Value = -3552.1426
$ws.Range("N42").Value = -6454
$ws.Range("H64").Value = 4500
$ws.Range("H67").Value = 4500
$ws.Range("H70").Value = 10173.2
$ws.Range("I70").Value = 8460
$ws.Range("J70").Value = 11315.333
$ws.Range("K70").Value = 25380
$ws.Range("L70").Value = 33945.999
$ws.Range("M70").Value = -25110
$ws.Range("N70").Value = -34485.999
$ws.Range("H73").Value = 10173.2
$ws.Range("I73").Value = 8460
$ws.Range("J73").Value = 11315.333
$ws.Range("K73").Value = 25380
$ws.Range("L73").Value = 33945.999
$ws.Range("M73").Value = -24444
$ws.Range("N73").Value = -35817.999
$ws.Range("H86").Value = 5558.952
$ws.Range("J86").Value = 6679.2856
$ws.Range("L86").Value = 6679.2856
$ws.Range("N86").Value = -8925.285599999999
$ws.Range("H89").Value = 5558.952
$ws.Range("J89").Value = 6679.2856
$ws.Range("L89").Value = 33396.428
$ws.Range("N89").Value = -44628.428
$ws.Range("H106").Value = 2096
$ws.Range("I106").Value = 2096
$ws.Range("K106").Value = 2096
$ws.Range("M106").Value = -1465
$ws.Range("H112").Value = 884.7273
$ws.Range("J112").Value = 884.7273
$ws.Range("L112").Value = 2654.1819
$ws.Range("N112").Value = -4870.1819
$ws.Range("H135").Value = 1456.6666
$ws.Range("I135").Value = 1522.1428
$ws.Range("J135").Value = 1227.5
$ws.Range("K135").Value = 13699.2852
$ws.Range("L135").Value = 11047.5
$ws.Range("M135").Value = -11164.2852
$ws.Range("N135").Value = -16117.5
$ws.Range("H137").Value = 1268.6
$ws.Range("I137").Value = 615.6667
$ws.Range("K137").Value = 1847.0001
$ws.Range("M137").Value = 702.9999
$ws.Range("H138").Value = 3460.3372
$ws.Range("I138").Value = 3818.8333
$ws.Range("J138").Value = 3365.4412
$ws.Range("K138").Value = 11456.4999
$ws.Range("L138").Value = 10096.3236
$ws.Range("M138").Value = -6316.499899999999
$ws.Range("N138").Value = -20376.3236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1801.5
$ws.Range("I132").Value = 1801.5
$ws.Range("K132").Value = 5404.5
$ws.Range("M132").Value = -2874.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39499.332
$ws.Range("I35").Value = 39499.332
$ws.Range("K35").Value = 39499.332
$ws.Range("M35").Value = -39189.332
$ws.Range("H94").Value = 526.1
$ws.Range("I94").Value = 507.875
$ws.Range("K94").Value = 507.875
$ws.Range("M94").Value = -56.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6384.778
$ws.Range("I31").Value = 6473
$ws.Range("K31").Value = 6473
$ws.Range("M31").Value = -6178
$ws.Range("H34").Value = 6384.778
$ws.Range("I34").Value = 6473
$ws.Range("K34").Value = 6473
$ws.Range("M34").Value = -6271
$ws.Range("H99").Value = 9934.111000000001
$ws.Range("I99").Value = 11195.728
$ws.Range("K99").Value = 11195.728
$ws.Range("M99").Value = -9697.727999999999
$ws.Range("H105").Value = 2864.1724
$ws.Range("I105").Value = 1989.4286
$ws.Range("K105").Value = 1989.4286
$ws.Range("M105").Value = -242.4286
$ws.Range("H122").Value = 2187.1428
$ws.Range("I122").Value = 1899.4
$ws.Range("K122").Value = 5698.200000000001
$ws.Range("M122").Value = -3248.200000000001
$ws.Range("H126").Value = 9934.111000000001
$ws.Range("I126").Value = 11195.728
$ws.Range("K126").Value = 33587.18399999999
$ws.Range("M126").Value = -31117.18399999999
$ws.Range("H132").Value = 1639.4736
$ws.Range("I132").Value = 1193.2727
$ws.Range("J132").Value = 2253
$ws.Range("K132").Value = 3579.8181
$ws.Range("L132").Value = 6759
$ws.Range("M132").Value = -1049.8181
$ws.Range("N132").Value = -11819
$ws.Range("H134").Value = 2380.0527
$ws.Range("I134").Value = 2295.1333
$ws.Range("K134").Value = 6885.3999
$ws.Range("M134").Value = -4350.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 153.33333
$ws.Range("I46").Value = 153.33333
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 459.99999
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -368.99999
$ws.Range("H88").Value = 13217.294
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 14978.143
$ws.Range("K88").Value = 15000
$ws.Range("L88").Value = 44934.429
$ws.Range("M88").Value = -14572
$ws.Range("N88").Value = -45790.429
$ws.Range("H91").Value = 13217.294
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 14978.143
$ws.Range("K91").Value = 15000
$ws.Range("L91").Value = 44934.429
$ws.Range("M91").Value = -13518
$ws.Range("N91").Value = -47898.429
$ws.Range("H122").Value = 666.1667
$ws.Range("J122").Value = 643.6667
$ws.Range("L122").Value = 5793.0003
$ws.Range("N122").Value = -10693.0003
$ws.Range("H134").Value = 16856.5
$ws.Range("J134").Value = 18984
$ws.Range("L134").Value = 56952
$ws.Range("N134").Value = -67092
$ws.Range("H137").Value = 4998
$ws.Range("J137").Value = 5044
$ws.Range("L137").Value = 15132
$ws.Range("N137").Value = -25332
$ws.Range("H140").Value = 5354.933
$ws.Range("I140").Value = 2302.182
$ws.Range("K140").Value = 6906.545999999999
$ws.Range("M140").Value = -1726.545999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 23380.5
$ws.Range("J92").Value = 23380.5
$ws.Range("L92").Value = 23380.5
$ws.Range("N92").Value = -27124.5
$ws.Range("H132").Value = 2222.5386
$ws.Range("I132").Value = 2099.4443
$ws.Range("K132").Value = 6298.3329
$ws.Range("M132").Value = -3768.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2746.5833
$ws.Range("I22").Value = 3114.6
$ws.Range("J22").Value = 2483.7144
$ws.Range("K22").Value = 3114.6
$ws.Range("L22").Value = 2483.7144
$ws.Range("M22").Value = -2819.6
$ws.Range("N22").Value = -3073.7144
$ws.Range("H27").Value = 2746.5833
$ws.Range("I27").Value = 3114.6
$ws.Range("J27").Value = 2483.7144
$ws.Range("K27").Value = 3114.6
$ws.Range("L27").Value = 2483.7144
$ws.Range("M27").Value = -3007.6
$ws.Range("N27").Value = -2697.7144
$ws.Range("H40").Value = 3576.682
$ws.Range("I40").Value = 3139.6155
$ws.Range("K40").Value = 3139.6155
$ws.Range("M40").Value = -3003.6155
$ws.Range("H45").Value = 29900
$ws.Range("I45").Value = 29900
$ws.Range("K45").Value = 29900
$ws.Range("M45").Value = -29493
$ws.Range("H100").Value = 5224.5
$ws.Range("J100").Value = 5225
$ws.Range("L100").Value = 5225
$ws.Range("N100").Value = -6307
$ws.Range("H104").Value = 18123
$ws.Range("J104").Value = 18123
$ws.Range("L104").Value = 18123
$ws.Range("N104").Value = -25111
$ws.Range("H122").Value = 6161.864
$ws.Range("I122").Value = 5466.3335
$ws.Range("J122").Value = 6996.5
$ws.Range("K122").Value = 16399.0005
$ws.Range("L122").Value = 20989.5
$ws.Range("M122").Value = -13949.0005
$ws.Range("N122").Value = -25889.5
$ws.Range("H132").Value = 3445.7334
$ws.Range("I132").Value = 2664
$ws.Range("J132").Value = 4339.143
$ws.Range("K132").Value = 7992
$ws.Range("L132").Value = 13017.429
$ws.Range("M132").Value = -5462
$ws.Range("N132").Value = -18077.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1699.3334
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 4684
$ws.Range("I126").Value = 2741.1667
$ws.Range("K126").Value = 8223.500100000001
$ws.Range("M126").Value = -5753.500100000001
$ws.Range("H132").Value = 3242.2354
$ws.Range("I132").Value = 3242.2354
$ws.Range("K132").Value = 9726.706200000001
$ws.Range("M132").Value = -7196.706200000001
$ws.Range("H136").Value = 2547
$ws.Range("I136").Value = 2547
$ws.Range("K136").Value = 7641
$ws.Range("M136").Value = -5091
